$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column D (Price) as Text so that values such as "2.68" or
# "221.15" are stored as strings (matching the source inlineStr cells)
# instead of being auto-converted to numbers by the Excel input parser.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '95.265.87'
$ws.Range("E2").Value = '  -1.92%  '

# Row 3
$ws.Range("D3").Value = '3.590.82'
$ws.Range("E3").Value = '  -3.16%  '

# Row 4
$ws.Range("D4").Value = '2.68'
$ws.Range("E4").Value = '  +31.36%  '

# Row 5
$ws.Range("E5").Value = '  +0.06%  '

# Row 6
$ws.Range("D6").Value = '221.15'
$ws.Range("E6").Value = '  -6.07%  '

# Row 7
$ws.Range("D7").Value = '634.97'
$ws.Range("E7").Value = '  -3.29%  '

# Row 8
$ws.Range("D8").Value = '0.419'
$ws.Range("E8").Value = '  -2.92%  '

# Row 9
$ws.Range("D9").Value = '1.18'
$ws.Range("E9").Value = '  +9.21%  '

# Row 11
$ws.Range("D11").Value = '3.584.46'
$ws.Range("E11").Value = '  -3.27%  '

# Row 12
$ws.Range("D12").Value = '47.79'
$ws.Range("E12").Value = '  +6.30%  '

# Row 13
$ws.Range("D13").Value = '0.212'
$ws.Range("E13").Value = '  +2.72%  '

# Row 14
$ws.Range("D14").Value = '0.0000291'
$ws.Range("E14").Value = '  -5.97%  '

# Row 15
$ws.Range("D15").Value = '6.45'
$ws.Range("E15").Value = '  -6.30%  '

# Row 16
$ws.Range("D16").Value = '4.261.60'
$ws.Range("E16").Value = '  -3.26%  '

# Row 17
$ws.Range("D17").Value = '94.964.53'
$ws.Range("E17").Value = '  -2.02%  '

# Row 18
$ws.Range("D18").Value = '22.79'
$ws.Range("E18").Value = '  +21.75%  '

# Row 19
$ws.Range("D19").Value = '8.85'
$ws.Range("E19").Value = '  -3.01%  '

# Row 20
$ws.Range("D20").Value = '13.72'
$ws.Range("E20").Value = '  +5.47%  '

# Row 21
$ws.Range("D21").Value = '3.584.52'
$ws.Range("E21").Value = '  -3.29%  '

# Row 22
$ws.Range("D22").Value = '0.288'
$ws.Range("E22").Value = '  +46.83%  '

# Row 23
$ws.Range("D23").Value = '0.535'
$ws.Range("E23").Value = '  +2.62%  '

# Row 24
$ws.Range("D24").Value = '509.53'
$ws.Range("E24").Value = '  -2.89%  '

# Row 25
$ws.Range("B25").Value = 'SuiNetwork'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D25").Value = '3.22'
$ws.Range("E25").Value = '  -6.82%  '

# Row 26
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").Value = '125.16'
$ws.Range("E26").Value = '  +16.79%  '

# Row 27
$ws.Range("D27").Value = '0.0000200'
$ws.Range("E27").Value = '  -10.46%  '

# Row 28
$ws.Range("D28").Value = '6.77'
$ws.Range("E28").Value = '  -1.69%  '

# Row 29
$ws.Range("D29").Value = '3.760.36'
$ws.Range("E29").Value = '  -3.81%  '

# Row 30
$ws.Range("D30").Value = '12.64'
$ws.Range("E30").Value = '  -6.27%  '

# Row 31
$ws.Range("D31").Value = '12.85'
$ws.Range("E31").Value = '  +1.83%  '

# Row 32
$ws.Range("D32").Value = '3.04'
$ws.Range("E32").Value = '  +0.34%  '

# Row 33
$ws.Range("E33").Value = '  +0.20%  '

# Row 34
$ws.Range("D34").Value = '0.617'
$ws.Range("E34").Value = '  +3.90%  '

# Row 35
$ws.Range("D35").Value = '0.180'

# Row 36
$ws.Range("B36").Value = 'Binance-PegBSC-USD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  -0.05%  '

# Row 37
$ws.Range("B37").Value = 'EthereumClassic'
$ws.Range("C37").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D37").Value = '32.51'
$ws.Range("E37").Value = '  -0.72%  '

# Row 38
$ws.Range("D38").Value = '1.75'
$ws.Range("E38").Value = '  -4.23%  '

# Row 39
$ws.Range("B39").Value = 'USDe'
$ws.Range("C39").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  +0.00%  '

# Row 40
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = '0.529'
$ws.Range("E40").Value = '  +5.90%  '

# Row 41
$ws.Range("D41").Value = '7.07'
$ws.Range("E41").Value = '  +3.99%  '

# Row 42
$ws.Range("D42").Value = '8.27'

# Row 43
$ws.Range("D43").Value = '578.40'
$ws.Range("E43").Value = '  -10.05%  '

# Row 44
$ws.Range("D44").Value = '0.0517'
$ws.Range("E44").Value = '  +12.90%  '

# Row 45
$ws.Range("D45").Value = '41.71'
$ws.Range("E45").Value = '  +4.17%  '

# Row 46
$ws.Range("D46").Value = '0.155'
$ws.Range("E46").Value = '  -6.61%  '

# Row 47
$ws.Range("D47").Value = '0.953'
$ws.Range("E47").Value = '  -1.02%  '

# Row 48
$ws.Range("D48").Value = '1.92'
$ws.Range("E48").Value = '  -5.09%  '

# Row 49
$ws.Range("D49").Value = '9.02'
$ws.Range("E49").Value = '  +3.33%  '

# Row 50
$ws.Range("D50").Value = '231.42'
$ws.Range("E50").Value = '  +11.80%  '

# Row 51
$ws.Range("D51").Value = '23.47'
$ws.Range("E51").Value = '  -0.55%  '

# Restore the default (Normal) style on column D now that the text values
# are committed, so no stray per-cell style survives the edit.
$priceCol.Style = "Normal"

